$d = $word.ActiveDocument

# Turn on revision tracking just for the "." insertion so the new text
# lands in its own run (matching a later, separate edit) instead of being
# silently coalesced into the existing run when the package is saved.
$d.TrackRevisions = $true
$para1 = $d.Paragraphs.Item(1)
$r1 = $para1.Range
$r1.Collapse(0)
$r1.InsertAfter(".")
$d.TrackRevisions = $false
$d.Revisions.AcceptAll() | Out-Null

# Add the new paragraph with its own sentence after the (now two-run)
# first paragraph.
$para1 = $d.Paragraphs.Item(1)
$r2 = $para1.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$r3 = $d.Paragraphs.Item(2).Range
$r3.Collapse(0)
$r3.InsertAfter("Leaflet-Part-2 not attempted in this challenge.")
